$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044871575449428
$ws.Cells.Item(2, 4).Value = 1.045513116010696
$ws.Cells.Item(2, 5).Value = 1.058222930568204
$ws.Cells.Item(2, 6).Value = 1.065153902051089
$ws.Cells.Item(2, 9).Value = 1.043521990393978
$ws.Cells.Item(2, 10).Value = 1.049934340087536
$ws.Cells.Item(2, 11).Value = 1.048281145596202
$ws.Cells.Item(2, 12).Value = 1.060955752903116
$ws.Cells.Item(2, 13).Value = 1.067867908212554
$ws.Cells.Item(2, 14).Value = 1.051425367279613

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04585432389959
$ws.Cells.Item(3, 4).Value = 1.046244273953502
$ws.Cells.Item(3, 5).Value = 1.059227385963914
$ws.Cells.Item(3, 6).Value = 1.066299642683366
$ws.Cells.Item(3, 9).Value = 1.043798993715062
$ws.Cells.Item(3, 10).Value = 1.050564355742953
$ws.Cells.Item(3, 11).Value = 1.048823968273799
$ws.Cells.Item(3, 12).Value = 1.061773715953557
$ws.Cells.Item(3, 13).Value = 1.068828166475832
$ws.Cells.Item(3, 14).Value = 1.052056277629524

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.046490428242637
$ws.Cells.Item(4, 4).Value = 1.046717447159297
$ws.Cells.Item(4, 5).Value = 1.05987794047193
$ws.Cells.Item(4, 6).Value = 1.067041926222849
$ws.Cells.Item(4, 9).Value = 1.043976953977678
$ws.Cells.Item(4, 10).Value = 1.050971601249028
$ws.Cells.Item(4, 11).Value = 1.049174607287866
$ws.Cells.Item(4, 12).Value = 1.062302983137035
$ws.Cells.Item(4, 13).Value = 1.069449843228557
$ws.Cells.Item(4, 14).Value = 1.05246410147093

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.046757894086728
$ws.Cells.Item(5, 4).Value = 1.046916383853242
$ws.Cells.Item(5, 5).Value = 1.060151577732052
$ws.Cells.Item(5, 6).Value = 1.067354200289691
$ws.Cells.Item(5, 9).Value = 1.044051461827065
$ws.Cells.Item(5, 10).Value = 1.0511427069807
$ws.Cells.Item(5, 11).Value = 1.049321870862598
$ws.Cells.Item(5, 12).Value = 1.06252548423294
$ws.Cells.Item(5, 13).Value = 1.069711273574512
$ws.Cells.Item(5, 14).Value = 1.052635450192366

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.046802805557525
$ws.Cells.Item(6, 4).Value = 1.046949787033701
$ws.Cells.Item(6, 5).Value = 1.060197531066808
$ws.Cells.Item(6, 6).Value = 1.06740664525432
$ws.Cells.Item(6, 9).Value = 1.044063954030816
$ws.Cells.Item(6, 10).Value = 1.051171430490799
$ws.Cells.Item(6, 11).Value = 1.049346588541322
$ws.Cells.Item(6, 12).Value = 1.062562842964852
$ws.Cells.Item(6, 13).Value = 1.069755173407336
$ws.Cells.Item(6, 14).Value = 1.052664214493144

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046494001947869
$ws.Cells.Item(7, 4).Value = 1.046720105304612
$ws.Cells.Item(7, 5).Value = 1.059881596260775
$ws.Cells.Item(7, 6).Value = 1.067046097987908
$ws.Cells.Item(7, 9).Value = 1.043977950760552
$ws.Cells.Item(7, 10).Value = 1.050973887966853
$ws.Cells.Item(7, 11).Value = 1.049176575601692
$ws.Cells.Item(7, 12).Value = 1.06230595621934
$ws.Cells.Item(7, 13).Value = 1.06945333616994
$ws.Cells.Item(7, 14).Value = 1.052466391436155

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04520365816177
$ws.Cells.Item(8, 4).Value = 1.045760200227594
$ws.Cells.Item(8, 5).Value = 1.058562265390425
$ws.Cells.Item(8, 6).Value = 1.065540921205862
$ws.Cells.Item(8, 9).Value = 1.043615869661435
$ws.Cells.Item(8, 10).Value = 1.050147342813375
$ws.Cells.Item(8, 11).Value = 1.048464719686556
$ws.Cells.Item(8, 12).Value = 1.061232189016526
$ws.Cells.Item(8, 13).Value = 1.068192364317593
$ws.Cells.Item(8, 14).Value = 1.051638672493753

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042931459343919
$ws.Cells.Item(9, 4).Value = 1.044069263712775
$ws.Cells.Item(9, 5).Value = 1.056242094335347
$ws.Cells.Item(9, 6).Value = 1.062895617244296
$ws.Cells.Item(9, 9).Value = 1.042968053197926
$ws.Cells.Item(9, 10).Value = 1.048687695073564
$ws.Cells.Item(9, 11).Value = 1.04720574126911
$ws.Cells.Item(9, 12).Value = 1.059340023560351
$ws.Cells.Item(9, 13).Value = 1.065972883996105
$ws.Cells.Item(9, 14).Value = 1.050176951886729

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.041417724241396
$ws.Cells.Item(10, 4).Value = 1.042942387756301
$ws.Cells.Item(10, 5).Value = 1.054698480414818
$ws.Cells.Item(10, 6).Value = 1.061136817965907
$ws.Cells.Item(10, 9).Value = 1.042529617735946
$ws.Cells.Item(10, 10).Value = 1.047712494041499
$ws.Cells.Item(10, 11).Value = 1.046363364004099
$ws.Cells.Item(10, 12).Value = 1.058078575349395
$ws.Cells.Item(10, 13).Value = 1.06449494592316
$ws.Cells.Item(10, 14).Value = 1.049200365957341

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040762516245141
$ws.Cells.Item(11, 4).Value = 1.04245454717273
$ws.Cells.Item(11, 5).Value = 1.054030836314732
$ws.Cells.Item(11, 6).Value = 1.060376366141589
$ws.Cells.Item(11, 9).Value = 1.042338218806388
$ws.Cells.Item(11, 10).Value = 1.047289728055492
$ws.Cells.Item(11, 11).Value = 1.045997886426492
$ws.Cells.Item(11, 12).Value = 1.057532358854647
$ws.Cells.Item(11, 13).Value = 1.063855394043797
$ws.Cells.Item(11, 14).Value = 1.048776999595142

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.04051918047335
$ws.Cells.Item(12, 4).Value = 1.0422733579165
$ws.Cells.Item(12, 5).Value = 1.053782956844091
$ws.Cells.Item(12, 6).Value = 1.06009406896952
$ws.Cells.Item(12, 9).Value = 1.042266891589695
$ws.Cells.Item(12, 10).Value = 1.047132619579287
$ws.Cells.Item(12, 11).Value = 1.045862023672068
$ws.Cells.Item(12, 12).Value = 1.05732946993861
$ws.Cells.Item(12, 13).Value = 1.063617897151986
$ws.Cells.Item(12, 14).Value = 1.04861966800688

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.040571375115414
$ws.Cells.Item(13, 4).Value = 1.042312222869267
$ws.Cells.Item(13, 5).Value = 1.053836122698672
$ws.Cells.Item(13, 6).Value = 1.06015461504247
$ws.Cells.Item(13, 9).Value = 1.042282202065321
$ws.Cells.Item(13, 10).Value = 1.04716632323984
$ws.Cells.Item(13, 11).Value = 1.045891171578824
$ws.Cells.Item(13, 12).Value = 1.057372990284259
$ws.Cells.Item(13, 13).Value = 1.063668838269418
$ws.Cells.Item(13, 14).Value = 1.048653419530497

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.040742401265105
$ws.Cells.Item(14, 4).Value = 1.042439569671869
$ws.Cells.Item(14, 5).Value = 1.054010344206926
$ws.Cells.Item(14, 6).Value = 1.060353027930601
$ws.Cells.Item(14, 9).Value = 1.042332327630717
$ws.Cells.Item(14, 10).Value = 1.047276742934672
$ws.Cells.Item(14, 11).Value = 1.045986658171546
$ws.Cells.Item(14, 12).Value = 1.057515587987086
$ws.Cells.Item(14, 13).Value = 1.063835761227417
$ws.Cells.Item(14, 14).Value = 1.048763996033961

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040847781115636
$ws.Cells.Item(15, 4).Value = 1.042518034444973
$ws.Cells.Item(15, 5).Value = 1.054117702851995
$ws.Cells.Item(15, 6).Value = 1.060475298994058
$ws.Cells.Item(15, 9).Value = 1.042363180761461
$ws.Cells.Item(15, 10).Value = 1.047344766310202
$ws.Cells.Item(15, 11).Value = 1.046045476301112
$ws.Cells.Item(15, 12).Value = 1.057603447177255
$ws.Cells.Item(15, 13).Value = 1.063938616083193
$ws.Cells.Item(15, 14).Value = 1.048832116010487

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.041461213507329
$ws.Cells.Item(16, 4).Value = 1.042974766392613
$ws.Cells.Item(16, 5).Value = 1.054742805614895
$ws.Cells.Item(16, 6).Value = 1.061187310264343
$ws.Cells.Item(16, 9).Value = 1.042542287544697
$ws.Cells.Item(16, 10).Value = 1.047740541178125
$ws.Cells.Item(16, 11).Value = 1.046387604361058
$ws.Cells.Item(16, 12).Value = 1.058114825948108
$ws.Cells.Item(16, 13).Value = 1.064537399434457
$ws.Cells.Item(16, 14).Value = 1.049228452924118

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041846070321198
$ws.Cells.Item(17, 4).Value = 1.043261290731744
$ws.Cells.Item(17, 5).Value = 1.055135117427229
$ws.Cells.Item(17, 6).Value = 1.061634236117174
$ws.Cells.Item(17, 9).Value = 1.042654220770453
$ws.Cells.Item(17, 10).Value = 1.047988667367835
$ws.Cells.Item(17, 11).Value = 1.046602019103319
$ws.Cells.Item(17, 12).Value = 1.058435600311516
$ws.Cells.Item(17, 13).Value = 1.064913109065454
$ws.Cells.Item(17, 14).Value = 1.049476931481478

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.042070574829105
$ws.Cells.Item(18, 4).Value = 1.043428425524992
$ws.Cells.Item(18, 5).Value = 1.055364018691398
$ws.Cells.Item(18, 6).Value = 1.061895028440162
$ws.Cells.Item(18, 9).Value = 1.042719359621278
$ws.Cells.Item(18, 10).Value = 1.048133347148944
$ws.Cells.Item(18, 11).Value = 1.046727013837308
$ws.Cells.Item(18, 12).Value = 1.058622702439698
$ws.Cells.Item(18, 13).Value = 1.06513229319898
$ws.Cells.Item(18, 14).Value = 1.049621816724471

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.042147129136157
$ws.Cells.Item(19, 4).Value = 1.043485415840758
$ws.Cells.Item(19, 5).Value = 1.055442080377146
$ws.Cells.Item(19, 6).Value = 1.061983970199878
$ws.Cells.Item(19, 9).Value = 1.042741544851505
$ws.Cells.Item(19, 10).Value = 1.048182671057654
$ws.Cells.Item(19, 11).Value = 1.046769621965283
$ws.Cells.Item(19, 12).Value = 1.058686499396675
$ws.Cells.Item(19, 13).Value = 1.065207035940945
$ws.Cells.Item(19, 14).Value = 1.049671210678789

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041804776339267
$ws.Cells.Item(20, 4).Value = 1.043230548342171
$ws.Cells.Item(20, 5).Value = 1.055093018588287
$ws.Cells.Item(20, 6).Value = 1.061586274026443
$ws.Cells.Item(20, 9).Value = 1.042642226898217
$ws.Cells.Item(20, 10).Value = 1.047962050740822
$ws.Cells.Item(20, 11).Value = 1.046579021642633
$ws.Cells.Item(20, 12).Value = 1.058401184268861
$ws.Cells.Item(20, 13).Value = 1.064872794935106
$ws.Cells.Item(20, 14).Value = 1.049450277055802

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.040692037285771
$ws.Cells.Item(21, 4).Value = 1.042402068751716
$ws.Cells.Item(21, 5).Value = 1.053959037184462
$ws.Cells.Item(21, 6).Value = 1.060294595627079
$ws.Cells.Item(21, 9).Value = 1.042317573334542
$ws.Cells.Item(21, 10).Value = 1.047244229136479
$ws.Cells.Item(21, 11).Value = 1.045958542730554
$ws.Cells.Item(21, 12).Value = 1.057473596508926
$ws.Cells.Item(21, 13).Value = 1.063786604886231
$ws.Cells.Item(21, 14).Value = 1.048731436062445

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.039992632017158
$ws.Cells.Item(22, 4).Value = 1.041881266370436
$ws.Cells.Item(22, 5).Value = 1.053246713671267
$ws.Cells.Item(22, 6).Value = 1.059483441133856
$ws.Cells.Item(22, 9).Value = 1.042112102052886
$ws.Cells.Item(22, 10).Value = 1.046792476137357
$ws.Cells.Item(22, 11).Value = 1.045567797754812
$ws.Cells.Item(22, 12).Value = 1.056890386673105
$ws.Cells.Item(22, 13).Value = 1.063104028049064
$ws.Cells.Item(22, 14).Value = 1.048279041522247

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.040363379258114
$ws.Cells.Item(23, 4).Value = 1.042157344232417
$ws.Cells.Item(23, 5).Value = 1.053624267545877
$ws.Cells.Item(23, 6).Value = 1.059913356977682
$ws.Cells.Item(23, 9).Value = 1.042221154020655
$ws.Cells.Item(23, 10).Value = 1.047031999588033
$ws.Cells.Item(23, 11).Value = 1.045774998172066
$ws.Cells.Item(23, 12).Value = 1.057199556982836
$ws.Cells.Item(23, 13).Value = 1.063465841294255
$ws.Cells.Item(23, 14).Value = 1.048518905123697

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041823435237853
$ws.Cells.Item(24, 4).Value = 1.043244439473886
$ws.Cells.Item(24, 5).Value = 1.055112041017667
$ws.Cells.Item(24, 6).Value = 1.061607945696839
$ws.Cells.Item(24, 9).Value = 1.042647646875976
$ws.Cells.Item(24, 10).Value = 1.047974077797235
$ws.Cells.Item(24, 11).Value = 1.046589413420732
$ws.Cells.Item(24, 12).Value = 1.058416735397471
$ws.Cells.Item(24, 13).Value = 1.064891011034539
$ws.Cells.Item(24, 14).Value = 1.049462321192014

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043518691177926
$ws.Cells.Item(25, 4).Value = 1.044506342033138
$ws.Cells.Item(25, 5).Value = 1.056841358544848
$ws.Cells.Item(25, 6).Value = 1.063578658596675
$ws.Cells.Item(25, 9).Value = 1.043136686605787
$ws.Cells.Item(25, 10).Value = 1.049065421679083
$ws.Cells.Item(25, 11).Value = 1.047531758810174
$ws.Cells.Item(25, 12).Value = 1.059829196646888
$ws.Cells.Item(25, 13).Value = 1.066546372390458
$ws.Cells.Item(25, 14).Value = 1.050555214907354
